$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

# Swap the shared-string content of row 24 and row 25 so that the
# "LOT2038" requirement now appears first (row 24) and the "LOT2028"
# requirement appears second (row 25).
$ws.Range("B24").Value = $lot2038
$ws.Range("C24").Value = $lot2038

$ws.Range("B25").Value = $lot2028
$ws.Range("C25").Value = $lot2028
